# "Generate Report for Handback"
# The handback file bb7ca973-2619-41a2-a4d2-ed0759aca5bf.md has been handed
# back and is now in sync with en-US, so update the localization status
# report across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for zh-cn (E) and de-de (F) on the
# bb7ca973... row (row 3)
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn detail sheet: Status (C), Latest Handback DateTime (K), Error Detail (P)
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-12 17:00:30"
$zhcn.Range("P3").Value = ""

# de-de detail sheet: Status (C), Latest Handback DateTime (K), Error Detail (P)
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-12 17:00:40"
$dede.Range("P3").Value = ""
